{"js": "// Replace each two-digit-division expression with its updated value.\n// The document stores each problem as a standalone run of text like\n// \"73\u00f78=\" inside its own table cell paragraph; every value is unique in\n// the document, so an exact-text search-and-replace is unambiguous.\nconst replacements = [\n  [\"73\u00f78=\", \"97\u00f76=\"],\n  [\"64\u00f72=\", \"98\u00f72=\"],\n  [\"23\u00f78=\", \"31\u00f77=\"],\n  [\"86\u00f76=\", \"67\u00f72=\"],\n  [\"99\u00f78=\", \"32\u00f77=\"],\n  [\"56\u00f78=\", \"65\u00f79=\"],\n  [\"41\u00f73=\", \"93\u00f78=\"],\n  [\"88\u00f76=\", \"68\u00f79=\"],\n  [\"83\u00f74=\", \"13\u00f73=\"],\n  [\"20\u00f72=\", \"58\u00f72=\"],\n  [\"57\u00f76=\", \"58\u00f73=\"],\n  [\"69\u00f76=\", \"83\u00f79=\"],\n  [\"35\u00f72=\", \"69\u00f79=\"],\n  [\"52\u00f79=\", \"75\u00f72=\"],\n  [\"49\u00f75=\", \"47\u00f73=\"],\n  [\"28\u00f74=\", \"14\u00f74=\"],\n  [\"58\u00f76=\", \"90\u00f76=\"],\n  [\"53\u00f72=\", \"45\u00f73=\"],\n  [\"84\u00f79=\", \"60\u00f79=\"],\n  [\"77\u00f74=\", \"87\u00f74=\"],\n  [\"23\u00f79=\", \"82\u00f76=\"],\n  [\"30\u00f72=\", \"91\u00f73=\"],\n  [\"38\u00f75=\", \"46\u00f79=\"],\n  [\"64\u00f78=\", \"96\u00f77=\"],\n  [\"76\u00f73=\", \"14\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update each two-digit-division expression in the table to its new value.\n# Every expression text (e.g. \"73\u00f78=\") occurs exactly once in the document,\n# so a plain Find/Replace (no wildcards) for each pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"73\u00f78=\", \"97\u00f76=\"),\n    @(\"64\u00f72=\", \"98\u00f72=\"),\n    @(\"23\u00f78=\", \"31\u00f77=\"),\n    @(\"86\u00f76=\", \"67\u00f72=\"),\n    @(\"99\u00f78=\", \"32\u00f77=\"),\n    @(\"56\u00f78=\", \"65\u00f79=\"),\n    @(\"41\u00f73=\", \"93\u00f78=\"),\n    @(\"88\u00f76=\", \"68\u00f79=\"),\n    @(\"83\u00f74=\", \"13\u00f73=\"),\n    @(\"20\u00f72=\", \"58\u00f72=\"),\n    @(\"57\u00f76=\", \"58\u00f73=\"),\n    @(\"69\u00f76=\", \"83\u00f79=\"),\n    @(\"35\u00f72=\", \"69\u00f79=\"),\n    @(\"52\u00f79=\", \"75\u00f72=\"),\n    @(\"49\u00f75=\", \"47\u00f73=\"),\n    @(\"28\u00f74=\", \"14\u00f74=\"),\n    @(\"58\u00f76=\", \"90\u00f76=\"),\n    @(\"53\u00f72=\", \"45\u00f73=\"),\n    @(\"84\u00f79=\", \"60\u00f79=\"),\n    @(\"77\u00f74=\", \"87\u00f74=\"),\n    @(\"23\u00f79=\", \"82\u00f76=\"),\n    @(\"30\u00f72=\", \"91\u00f73=\"),\n    @(\"38\u00f75=\", \"46\u00f79=\"),\n    @(\"64\u00f78=\", \"96\u00f77=\"),\n    @(\"76\u00f73=\", \"14\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$find.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$find.Replacement.Text, [ref]2) | Out-Null\n}\n"}
